$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete data rows 2-5, leaving only the header row
$ws.Rows("2:5").Delete()

# Add new column header "execution_time" in Q1
$ws.Range("Q1").Value = "execution_time"
$ws.Range("B1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
